# Auto-generated edit script: refresh market-data columns (H-N) across the
# "Leve Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# the latest scheduled market-data pull. Values are literal numbers (no
# formulas on this sheet), so each touched cell is written directly.
#
# Cells whose new value is absent entirely in the target (diff shows the
# whole <c> element removed) are cleared instead of set to 0/blank text so
# the cell stops existing in the saved sheet, matching upstream behaviour
# for these zero-result derived columns.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 140.75
$ws.Range("I33").Value = 142.84616
$ws.Range("J33").Value = 131.66667
$ws.Range("K33").Value = 142.84616
$ws.Range("L33").Value = 131.66667
$ws.Range("M33").Value = 86.15384
$ws.Range("N33").Value = -589.6666700000001
$ws.Range("H40").Value = 1913.3529
$ws.Range("I40").Value = 1737.9412
$ws.Range("K40").Value = 1737.9412
$ws.Range("M40").Value = -1562.9412
$ws.Range("H74").Value = 121706.695
$ws.Range("I74").Value = 131432.25
$ws.Range("K74").Value = 131432.25
$ws.Range("M74").Value = -130496.25
$ws.Range("H77").Value = 121706.695
$ws.Range("I77").Value = 131432.25
$ws.Range("K77").Value = 657161.25
$ws.Range("M77").Value = -652481.25
$ws.Range("H86").Value = 1969
$ws.Range("I86").Value = 1561.3334
$ws.Range("J86").Value = 2274.75
$ws.Range("K86").Value = 1561.3334
$ws.Range("L86").Value = 2274.75
$ws.Range("M86").Value = -438.3334
$ws.Range("N86").Value = -4520.75
$ws.Range("H89").Value = 1969
$ws.Range("I89").Value = 1561.3334
$ws.Range("J89").Value = 2274.75
$ws.Range("K89").Value = 7806.666999999999
$ws.Range("L89").Value = 11373.75
$ws.Range("M89").Value = -2190.666999999999
$ws.Range("N89").Value = -22605.75
$ws.Range("H92").Value = 322.5
$ws.Range("I92").Value = 322.5
$ws.Range("K92").Value = 322.5
$ws.Range("M92").Value = 925.5
$ws.Range("H94").Value = 1632.6
$ws.Range("I94").Value = 1543.25
$ws.Range("K94").Value = 1543.25
$ws.Range("M94").Value = -1092.25
$ws.Range("H98").Value = 1637.5
$ws.Range("I98").Value = 1850
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1850
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -352
$ws.Range("N98").Value = -3996
$ws.Range("H100").Value = 1010.55554
$ws.Range("I100").Value = 1011.875
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1011.875
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -470.875
$ws.Range("N100").Value = -2082
$ws.Range("I101").Value = 25004722
$ws.Range("J101").Value = 995
$ws.Range("K101").Value = 75014166
$ws.Range("L101").Value = 2985
$ws.Range("M101").Value = -75012544
$ws.Range("N101").Value = -6229
$ws.Range("H107").Value = 437.66666
$ws.Range("I107").Value = 486.5
$ws.Range("K107").Value = 486.5
$ws.Range("M107").Value = 1433.5
$ws.Range("H122").Value = 1637.5
$ws.Range("I122").Value = 1850
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5550
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3100
$ws.Range("N122").Value = -7900
$ws.Range("H127").Value = 1551.3572
$ws.Range("I127").Value = 1616.25
$ws.Range("K127").Value = 4848.75
$ws.Range("M127").Value = 111.25
$ws.Range("H129").Value = 3371
$ws.Range("I129").Value = 3666
$ws.Range("J129").Value = 3302.923
$ws.Range("K129").Value = 10998
$ws.Range("L129").Value = 9908.769
$ws.Range("M129").Value = -5998
$ws.Range("N129").Value = -19908.769
$ws.Range("H137").Value = 5362.636
$ws.Range("I137").Value = 3713.7144
$ws.Range("J137").Value = 8248.25
$ws.Range("K137").Value = 11141.1432
$ws.Range("L137").Value = 24744.75
$ws.Range("M137").Value = -8591.143199999999
$ws.Range("N137").Value = -29844.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2491
$ws.Range("J88").Value = 2491
$ws.Range("L88").Value = 2491
$ws.Range("N88").Value = -3303
$ws.Range("H91").Value = 2491
$ws.Range("J91").Value = 2491
$ws.Range("L91").Value = 2491
$ws.Range("N91").Value = -5299

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 609
$ws.Range("I94").Value = 645.6667
$ws.Range("K94").Value = 645.6667
$ws.Range("M94").Value = -194.6667
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 98.3
$ws.Range("I7").Value = 120.666664
$ws.Range("K7").Value = 120.666664
$ws.Range("M7").Value = -7.666663999999997
$ws.Range("H16").Value = 3624.125
$ws.Range("I16").Value = 4996.5
$ws.Range("J16").Value = 3166.6667
$ws.Range("K16").Value = 4996.5
$ws.Range("L16").Value = 3166.6667
$ws.Range("M16").Value = -4709.5
$ws.Range("N16").Value = -3740.6667
$ws.Range("H93").Value = 99
$ws.Range("J93").Value = 99
$ws.Range("L93").Value = 99
$ws.Range("N93").Value = -3843
$ws.Range("H99").Value = 1538.2222
$ws.Range("J99").Value = 2025.5
$ws.Range("L99").Value = 2025.5
$ws.Range("N99").Value = -5021.5
$ws.Range("H105").Value = 3720
$ws.Range("I105").Value = 2550
$ws.Range("K105").Value = 2550
$ws.Range("M105").Value = -803
$ws.Range("H113").Value = 3624.125
$ws.Range("I113").Value = 4996.5
$ws.Range("J113").Value = 3166.6667
$ws.Range("K113").Value = 4996.5
$ws.Range("L113").Value = 3166.6667
$ws.Range("M113").Value = -2826.5
$ws.Range("N113").Value = -7506.6667
$ws.Range("H126").Value = 1538.2222
$ws.Range("J126").Value = 2025.5
$ws.Range("L126").Value = 6076.5
$ws.Range("N126").Value = -11016.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 788.5
$ws.Range("I69").Value = 789
$ws.Range("J69").Value = 788
$ws.Range("K69").Value = 2367
$ws.Range("L69").Value = 2364
$ws.Range("M69").Value = -1556
$ws.Range("N69").Value = -3986
$ws.Range("H72").Value = 788.5
$ws.Range("I72").Value = 789
$ws.Range("J72").Value = 788
$ws.Range("K72").Value = 7101
$ws.Range("L72").Value = 7092
$ws.Range("M72").Value = -3045
$ws.Range("N72").Value = -15204
$ws.Range("H81").Value = 2300
$ws.Range("J81").Value = 2400
$ws.Range("L81").Value = 7200
$ws.Range("N81").Value = -9446
$ws.Range("H84").Value = 2300
$ws.Range("J84").Value = 2400
$ws.Range("L84").Value = 21600
$ws.Range("N84").Value = -32832

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 42499.5
$ws.Range("J69").Value = 42499.5
$ws.Range("L69").Value = 42499.5
$ws.Range("N69").Value = -43997.5
$ws.Range("H70").Value = 4985
$ws.Range("J70").Value = 4985
$ws.Range("L70").Value = 4985
$ws.Range("N70").Value = -5525
$ws.Range("H72").Value = 42499.5
$ws.Range("J72").Value = 42499.5
$ws.Range("L72").Value = 127498.5
$ws.Range("N72").Value = -134986.5
$ws.Range("H73").Value = 4985
$ws.Range("J73").Value = 4985
$ws.Range("L73").Value = 4985
$ws.Range("N73").Value = -6857
$ws.Range("H98").Value = 14322.857
$ws.Range("J98").Value = 14322.857
$ws.Range("L98").Value = 14322.857
$ws.Range("N98").Value = -20312.857
$ws.Range("H102").Value = 1704.6428
$ws.Range("I102").Value = 1655.5
$ws.Range("J102").Value = 1999.5
$ws.Range("K102").Value = 1655.5
$ws.Range("L102").Value = 1999.5
$ws.Range("M102").Value = -33.5
$ws.Range("N102").Value = -5243.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 55500
$ws.Range("I74").Value = 55500
$ws.Range("K74").Value = 55500
$ws.Range("M74").Value = -54502
$ws.Range("H75").Value = 253586.5
$ws.Range("I75").Value = 7000
$ws.Range("J75").Value = 500173
$ws.Range("K75").Value = 7000
$ws.Range("L75").Value = 500173
$ws.Range("M75").Value = -6064
$ws.Range("N75").Value = -502045
$ws.Range("H77").Value = 55500
$ws.Range("I77").Value = 55500
$ws.Range("K77").Value = 166500
$ws.Range("M77").Value = -161508
$ws.Range("H78").Value = 253586.5
$ws.Range("I78").Value = 7000
$ws.Range("J78").Value = 500173
$ws.Range("K78").Value = 21000
$ws.Range("L78").Value = 1500519
$ws.Range("M78").Value = -16320
$ws.Range("N78").Value = -1509879
$ws.Range("H82").Value = 3470
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 3470
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H93").Value = 2524.1667
$ws.Range("I93").Value = 2249.2
$ws.Range("J93").Value = 3899
$ws.Range("K93").Value = 2249.2
$ws.Range("L93").Value = 3899
$ws.Range("M93").Value = -1001.2
$ws.Range("N93").Value = -6395
$ws.Range("H122").Value = 9610.333000000001
$ws.Range("I122").Value = 13839.6
$ws.Range("K122").Value = 41518.8
$ws.Range("M122").Value = -39068.8
$ws.Range("H136").Value = 1966.6666
$ws.Range("I136").Value = 1850.75
$ws.Range("J136").Value = 2198.5
$ws.Range("K136").Value = 5552.25
$ws.Range("L136").Value = 6595.5
$ws.Range("M136").Value = -3002.25
$ws.Range("N136").Value = -11695.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 11940
$ws.Range("I7").Value = 320
$ws.Range("J7").Value = 17750
$ws.Range("K7").Value = 320
$ws.Range("L7").Value = 17750
$ws.Range("M7").Value = -207
$ws.Range("N7").Value = -17976
$ws.Range("H74").Value = 19847.25
$ws.Range("I74").Value = 18477.5
$ws.Range("K74").Value = 18477.5
$ws.Range("M74").Value = -17541.5
$ws.Range("H77").Value = 19847.25
$ws.Range("I77").Value = 18477.5
$ws.Range("K77").Value = 55432.5
$ws.Range("M77").Value = -50752.5
$ws.Range("H96").Value = 1710.25
$ws.Range("I96").Value = 1936.2
$ws.Range("K96").Value = 1936.2
$ws.Range("M96").Value = -563.2
$ws.Range("H114").Value = 22455
$ws.Range("J114").Value = 22455
$ws.Range("L114").Value = 22455
$ws.Range("N114").Value = -31133

